$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'245.47"
# Row 4
$ws.Range("D4").Value = "'5.131"
# Row 5
$ws.Range("D5").Value = "'0.05584"
# Row 6
$ws.Range("D6").Value = "'6.482"
# Row 8
$ws.Range("D8").Value = "'0.8176"
# Row 9
$ws.Range("D9").Value = "'0.8408"
# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1336"
$ws.Range("E10").Value = "9WazirXWRX"
# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03210"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
# Row 12
$ws.Range("D12").Value = "'0.02859"
# Row 13
$ws.Range("D13").Value = "'0.09384"
# Row 14
$ws.Range("D14").Value = "'0.001510"
# Row 15
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005949"
$ws.Range("E15").Value = "14OneONE"
# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006177"
$ws.Range("E16").Value = "15TigerCashTCH"
# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.528"
$ws.Range("E17").Value = "16LEOLEO"
# Row 18
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.022"
$ws.Range("E18").Value = "17BTSETokenBTSE"
# Row 19
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3179"
$ws.Range("E19").Value = "18BitpandaEcosystemTokenBEST"
# Row 20
$ws.Range("B20").Value = "MandalaExchangeToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D20").Value = "'0.06970"
$ws.Range("E20").Value = "19MandalaExchangeTokenMDX"
# Row 22
$ws.Range("D22").Value = "'3.746"
# Row 23
$ws.Range("D23").Value = "'0.04711"
# Row 25
$ws.Range("D25").Value = "'0.001247"
# Row 26
$ws.Range("D26").Value = "'0.004609"
# Row 27
$ws.Range("D27").Value = "'0.00009695"
# Row 40
$ws.Range("D40").Value = "'0.03657"
# Row 41
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1368"
$ws.Range("E41").Value = "40BKEXTokenBKKBestin24h"
# Row 42
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002523"
$ws.Range("E42").Value = "41CEJICEJI"
# Row 43
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003378"
$ws.Range("E43").Value = "42KickTokenKICK"
# Row 44
$ws.Range("D44").Value = "'0.007648"
# Row 45
$ws.Range("D45").Value = "'0.00005316"
